$d = $word.ActiveDocument
$d.Bookmarks.ShowHidden = $true

# --- 1) Remove the stale _GoBack bookmark (it currently sits in the empty
#        first paragraph at the top of the document). ---
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2) Rewrite the "Update Medication Types" precondition paragraph. ---
$r = $d.Content
$r.Find.Execute("User must have loaded the medication type into the Medication Details form.")
$start = $r.Start
$r.Text = ""

$pos = $start
$parts = @("User must have ", "located", " the ", "medication", " type ", "via Search Medication Types")
foreach ($t in $parts) {
  $cur = $d.Range($pos, $pos)
  $cur.InsertAfter($t)
  $pos = $pos + $t.Length
}

# New _GoBack bookmark right after "via Search Medication Types" (zero width),
# matching the location left by the author's last edit.
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Trailing period that closes the sentence.
$tail = $d.Range($pos, $pos)
$tail.InsertAfter(".")

# --- 3) Delete the now-superfluous empty paragraph that followed
#        "User must know what the field is to be updated to." ---
$r2 = $d.Content
$r2.Find.Execute("User must know what the field is to be updated to.")
$para = $r2.Paragraphs(1)
$nextPara = $para.Next()
if ($nextPara.Range.Text -eq "`r") {
  $nextPara.Range.Delete()
}

Write-Output "done"
